$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.1101440972136105
    "C2" = 2.223899453410976
    "D2" = 0.1077033850415691
    "B3" = 0.1475133730083216
    "C3" = 3.717662391659502
    "D3" = 0.2533868249348961
    "B4" = 0.0782670147820092
    "C4" = 4.555398714672766
    "D4" = 0.1672945128040517
    "B5" = 0.1047504286283093
    "C5" = 6.262020263624557
    "D5" = 0.1192100483910533
    "B6" = 0.05231663667223318
    "C6" = 7.271091143955386
    "D6" = 0.1489664110061631
    "B7" = 0.128124222806019
    "C7" = 8.593461662935168
    "D7" = 0.2832036018017433
    "B8" = 0.0791793114883132
    "C8" = 10.02817605633033
    "D8" = 0.1042301865793677
    "B9" = 0.120806835183015
    "C9" = 11.12898899860832
    "D9" = 0.1658036729864698
    "B10" = 0.05962567344995266
    "C10" = 12.33068157201087
    "D10" = 0.2783023277974669
    "B11" = 0.08795022932940902
    "C11" = 13.64120324978665
    "D11" = 0.1144082282813273
    "B12" = 0.05218613861180674
    "C12" = 14.82355532846044
    "D12" = 0.2269846343965896
    "B13" = 0.06130342491592651
    "C13" = 15.67666829356109
    "D13" = 0.1988834766327157
    "B14" = 0.09267940946585032
    "C14" = 17.45045298248814
    "D14" = 0.2491243864510206
    "B15" = 0.133482360147976
    "C15" = 18.32610243922323
    "D15" = 0.2495338327241333
    "B16" = 0.1138516941377553
    "C16" = 19.57407349875002
    "D16" = 0.2722374951570674
    "B17" = 0.07539235467898789
    "C17" = 21.24512124444964
    "D17" = 0.2486234090018789
    "B18" = 0.05463635605916428
    "C18" = 22.22618526947914
    "D18" = 0.1669028779207616
    "B19" = 0.07971880613462326
    "C19" = 23.21027939671141
    "D19" = 0.2712061299813804
    "B20" = 0.1400059669297661
    "C20" = 24.9889010220153
    "D20" = 0.2446167853980206
    "B21" = 0.05120723743929776
    "C21" = 25.95710958700444
    "D21" = 0.2617788206634494
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
